$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    # Force the cell to be written as text so numeric-looking strings
    # (e.g. "207.82") are not coerced into floating point numbers,
    # then restore the default style so no stray formatting is left behind.
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 2 4 '27.466.47'
Set-TextValue 2 5 '  -0.03%  '
Set-TextValue 3 4 '1.568.41'
Set-TextValue 3 5 '  +0.08%  '
Set-TextValue 4 5 '  -0.55%  '
Set-TextValue 5 4 '207.82'
Set-TextValue 5 5 '  +1.81%  '
Set-TextValue 6 5 '  -0.30%  '
Set-TextValue 7 5 '  -0.48%  '
Set-TextValue 8 4 '21.98'
Set-TextValue 8 5 '  -0.16%  '
Set-TextValue 9 4 '0.248'
Set-TextValue 9 5 '  -1.04%  '
Set-TextValue 10 4 '0.0589'
Set-TextValue 10 5 '  +0.73%  '
Set-TextValue 11 5 '  +0.59%  '
Set-TextValue 12 4 '1.793.34'
Set-TextValue 12 5 '  -0.34%  '
Set-TextValue 13 4 '1.572.86'
Set-TextValue 13 5 '  +0.38%  '
Set-TextValue 14 5 '  +0.20%  '
Set-TextValue 15 5 '  -1.49%  '
Set-TextValue 16 4 '63.29'
Set-TextValue 16 5 '  +1.29%  '
Set-TextValue 17 4 '27.476.40'
Set-TextValue 17 5 '  +0.00%  '
Set-TextValue 18 4 '213.95'
Set-TextValue 18 5 '  -0.40%  '
Set-TextValue 19 4 '0.0₃0689'
Set-TextValue 19 5 '  +0.80%  '
Set-TextValue 20 5 '  +0.23%  '
Set-TextValue 21 5 '  -0.37%  '
Set-TextValue 22 5 '  +0.25%  '
Set-TextValue 23 4 '9.54'
Set-TextValue 23 5 '  +1.09%  '
Set-TextValue 24 5 '  +1.61%  '
Set-TextValue 25 4 '153.14'
Set-TextValue 25 5 '  +0.06%  '
Set-TextValue 26 4 '6.80'
Set-TextValue 26 5 '  +2.13%  '
Set-TextValue 27 5 '  -0.67%  '
Set-TextValue 28 5 '  +0.62%  '
Set-TextValue 29 5 '  -1.16%  '
Set-TextValue 30 4 '1.15'
Set-TextValue 30 5 '  +0.50%  '
Set-TextValue 31 5 '  +1.86%  '
Set-TextValue 32 5 '  -0.30%  '
Set-TextValue 33 4 '1.362.49'
Set-TextValue 33 5 '  +0.45%  '
Set-TextValue 34 5 '  +1.26%  '
Set-TextValue 35 5 '  +3.13%  '
Set-TextValue 36 5 '  +1.29%  '
Set-TextValue 37 5 '  -0.35%  '
Set-TextValue 38 4 '0.0167'
Set-TextValue 38 5 '  +2.16%  '
Set-TextValue 39 5 '  +0.16%  '
Set-TextValue 40 4 '0.821'
Set-TextValue 40 5 '  +2.72%  '
Set-TextValue 41 5 '  -0.40%  '
Set-TextValue 42 4 '0.973'
Set-TextValue 42 5 '  -0.10%  '
Set-TextValue 43 5 '  +2.25%  '
Set-TextValue 44 4 '64.15'
Set-TextValue 44 5 '  +1.99%  '
Set-TextValue 45 5 '  +0.75%  '
Set-TextValue 46 5 '  -2.60%  '
Set-TextValue 47 4 '1.705.24'
Set-TextValue 47 5 '  -0.66%  '
Set-TextValue 48 4 '85.51'
Set-TextValue 48 5 '  -1.03%  '
Set-TextValue 49 4 '0.0₇0989'
Set-TextValue 49 5 '  +2.57%  '
Set-TextValue 50 4 '0.0953'
Set-TextValue 50 5 '  -0.64%  '
Set-TextValue 51 5 '  -0.09%  '
